# Update column G ("K") values on Sheet1 to reflect regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 2
    13 = 3
    14 = 1
    16 = 1
    18 = 1
    20 = 1
    21 = 0
    22 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
